$d = $word.ActiveDocument

# The second paragraph currently holds 23 repetitions of "A paragraph,".
# Replace its whole range (between the existing spell-check proofErr
# markers and the paragraph mark) with the m:for field-code runs plus
# the bold red error message, mirroring what M2Doc generates when a
# repetition tag is missing its iteration variable.
$p = $d.Paragraphs(2)
$r = $d.Range($p.Range.Start, $p.Range.End)

$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:fldChar w:fldCharType="begin"/></w:r><w:r><w:instrText>m</w:instrText></w:r><w:r><w:instrText xml:space="preserve">:for </w:instrText></w:r><w:r><w:instrText>|</w:instrText></w:r><w:r><w:instrText xml:space="preserve"> </w:instrText></w:r><w:r><w:instrText>self.e</w:instrText></w:r><w:r><w:instrText>Classifiers</w:instrText></w:r><w:r><w:instrText xml:space="preserve"> </w:instrText></w:r><w:r><w:fldChar w:fldCharType="end"/></w:r><w:r><w:rPr><w:b w:val="true"/><w:color w:val="FF0000"/></w:rPr><w:t>Invalid for statement: Malformed tag m:for : no iteration variable specified.</w:t></w:r><w:r><w:fldChar w:fldCharType="begin"/></w:r><w:r><w:instrText xml:space="preserve"> </w:instrText></w:r><w:r><w:instrText>m</w:instrText></w:r><w:r><w:instrText xml:space="preserve">:endfor </w:instrText></w:r><w:r><w:fldChar w:fldCharType="end"/></w:r></w:p>
'@

$r.InsertXML($xml)
